# Natmi following Dr Hou advice
# ---------------------------------------------------------------
# The NATMI ligand-receptor analysis (Efna4 -> Epha1) was re-run with
# an additional "M2" cluster included among the *sending* clusters
# (it was already present as a target cluster). This grows the
# sender x target grid from 3x4 (13 data rows) to 4x4 (16 data rows),
# and refreshes every computed statistic (detection rates, average /
# total expression, specificity scores, edge weights, ...) for all
# sender/target cluster combinations.
#
# Below we just rewrite the whole data block (rows 2-17, columns A-T)
# with the refreshed values; the header row (row 1) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efna4", "Epha1", "ECs", 3, 1, 0.8339443333333333, 2.501833, 0.305203232930382, 0.3052032329303821, 3, 1, 3.746298666666667, 11.238896, 0.2944715732214049, 0.294471573221405, 3.124204544040889, 28.117840896368, 0.0898736761532685, 0.08987367615326854),
    @("ECs", "Efna4", "Epha1", "FAPs", 3, 1, 0.8339443333333333, 2.501833, 0.305203232930382, 0.3052032329303821, 3, 1, 4.370261666666667, 13.110785, 0.3435171466234404, 0.3435171466234404, 3.644554952100556, 32.800994568905, 0.1048425437164941, 0.1048425437164941),
    @("ECs", "Efna4", "Epha1", "M2", 3, 1, 0.8339443333333333, 2.501833, 0.305203232930382, 0.3052032329303821, 3, 1, 1.837969333333333, 5.513908, 0.1444705212467569, 0.1444705212467569, 1.532764110373778, 13.794876993364, 0.04409287014764764, 0.04409287014764765),
    @("ECs", "Efna4", "Epha1", "sCs", 3, 1, 0.8339443333333333, 2.501833, 0.305203232930382, 0.3052032329303821, 3, 1, 2.767576666666667, 8.30273, 0.2175407589083977, 0.2175407589083977, 2.308004878232222, 20.77204390409, 0.0663941429129718, 0.06639414291297181),
    @("FAPs", "Efna4", "Epha1", "ECs", 3, 1, 1.537538333333333, 4.612615, 0.5627014314157557, 0.5627014314157558, 3, 1, 3.746298666666667, 11.238896, 0.2944715732214049, 0.294471573221405, 5.760077808115555, 51.84070027304, 0.1656995757629341, 0.1656995757629341),
    @("FAPs", "Efna4", "Epha1", "FAPs", 3, 1, 1.537538333333333, 4.612615, 0.5627014314157557, 0.5627014314157558, 3, 1, 4.370261666666667, 13.110785, 0.3435171466234404, 0.3435171466234404, 6.719444839197222, 60.475003552775, 0.1932975901208659, 0.193297590120866),
    @("FAPs", "Efna4", "Epha1", "M2", 3, 1, 1.537538333333333, 4.612615, 0.5627014314157557, 0.5627014314157558, 3, 1, 1.837969333333333, 5.513908, 0.1444705212467569, 0.1444705212467569, 2.825948305491111, 25.43353474942, 0.08129376910293043, 0.08129376910293044),
    @("FAPs", "Efna4", "Epha1", "sCs", 3, 1, 1.537538333333333, 4.612615, 0.5627014314157557, 0.5627014314157558, 3, 1, 2.767576666666667, 8.30273, 0.2175407589083977, 0.2175407589083977, 4.255255215438889, 38.29729693895, 0.1224104964290252, 0.1224104964290252),
    @("M2", "Efna4", "Epha1", "ECs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01459937449899473, 0.01459937449899473, 3, 1, 3.746298666666667, 11.238896, 0.2944715732214049, 0.294471573221405, 0.1494460976444444, 1.3450148788, 0.004299100776767438, 0.00429910077676744),
    @("M2", "Efna4", "Epha1", "FAPs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01459937449899473, 0.01459937449899473, 3, 1, 4.370261666666667, 13.110785, 0.3435171466234404, 0.3435171466234404, 0.1743370216527778, 1.569033194875, 0.005015135470381689, 0.005015135470381689),
    @("M2", "Efna4", "Epha1", "M2", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01459937449899473, 0.01459937449899473, 3, 1, 1.837969333333333, 5.513908, 0.1444705212467569, 0.1444705212467569, 0.07331965998888888, 0.6598769399, 0.002109179243746378, 0.002109179243746378),
    @("M2", "Efna4", "Epha1", "sCs", 1, 0.3333333333333333, 0.03989166666666667, 0.119675, 0.01459937449899473, 0.01459937449899473, 3, 1, 2.767576666666667, 8.30273, 0.2175407589083977, 0.2175407589083977, 0.1104032458611111, 0.99362921275, 0.003175959008099222, 0.003175959008099223),
    @("sCs", "Efna4", "Epha1", "ECs", 3, 1, 0.3210486666666667, 0.9631460000000001, 0.1174959611548676, 0.1174959611548676, 3, 1, 3.746298666666667, 11.238896, 0.2944715732214049, 0.294471573221405, 1.202744191868445, 10.824697726816, 0.03459922052843493, 0.03459922052843495),
    @("sCs", "Efna4", "Epha1", "FAPs", 3, 1, 0.3210486666666667, 0.9631460000000001, 0.1174959611548676, 0.1174959611548676, 3, 1, 4.370261666666667, 13.110785, 0.3435171466234404, 0.3435171466234404, 1.403066681067778, 12.62760012961, 0.0403618773156987, 0.04036187731569871),
    @("sCs", "Efna4", "Epha1", "M2", 3, 1, 0.3210486666666667, 0.9631460000000001, 0.1174959611548676, 0.1174959611548676, 3, 1, 1.837969333333333, 5.513908, 0.1444705212467569, 0.1444705212467569, 0.5900776038408889, 5.310698434568, 0.01697470275243241, 0.01697470275243242),
    @("sCs", "Efna4", "Epha1", "sCs", 3, 1, 0.3210486666666667, 0.9631460000000001, 0.1174959611548676, 0.1174959611548676, 3, 1, 2.767576666666667, 8.30273, 0.2175407589083977, 0.2175407589083977, 0.8885267987311113, 7.996741188580001, 0.02556016055830151, 0.02556016055830152)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
